$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 379 (weekly update: two new
# price records for "Pepino ensalada" dated 2022-11-25, serial 44890),
# pushing the existing rows 379-386 down to 381-388.
$ws.Rows.Item(379).Insert()
$ws.Rows.Item(379).Insert()

# New row 379 - Calidad "Primera"
$ws.Cells.Item(379, 1).Value2 = 1
$ws.Cells.Item(379, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(379, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(379, 4).Value2 = 44890
$ws.Cells.Item(379, 5).Value2 = 15
$ws.Cells.Item(379, 6).Value2 = 100112043
$ws.Cells.Item(379, 7).Value2 = "Pepino ensalada"
$ws.Cells.Item(379, 8).Value2 = "Sin especificar"
$ws.Cells.Item(379, 9).Value2 = "Primera"
$ws.Cells.Item(379, 10).Value2 = 470
$ws.Cells.Item(379, 11).Value2 = 7000
$ws.Cells.Item(379, 12).Value2 = 8000
$ws.Cells.Item(379, 13).Value2 = 7468
$ws.Cells.Item(379, 14).Value2 = "$/caja 70 unidades"
$ws.Cells.Item(379, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(379, 16).Value2 = 107
$ws.Cells.Item(379, 17).Value2 = 70
$ws.Cells.Item(379, 18).Value2 = "Hortaliza"

# New row 380 - Calidad "Segunda"
$ws.Cells.Item(380, 1).Value2 = 1
$ws.Cells.Item(380, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(380, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(380, 4).Value2 = 44890
$ws.Cells.Item(380, 5).Value2 = 15
$ws.Cells.Item(380, 6).Value2 = 100112043
$ws.Cells.Item(380, 7).Value2 = "Pepino ensalada"
$ws.Cells.Item(380, 8).Value2 = "Sin especificar"
$ws.Cells.Item(380, 9).Value2 = "Segunda"
$ws.Cells.Item(380, 10).Value2 = 370
$ws.Cells.Item(380, 11).Value2 = 4000
$ws.Cells.Item(380, 12).Value2 = 5000
$ws.Cells.Item(380, 13).Value2 = 4486
$ws.Cells.Item(380, 14).Value2 = "$/caja 100 unidades"
$ws.Cells.Item(380, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(380, 16).Value2 = 45
$ws.Cells.Item(380, 17).Value2 = 100
$ws.Cells.Item(380, 18).Value2 = "Hortaliza"
